# Fix footnote table symbols from symbols (*, †, ‡) to superscript letters (ᵃ, ᵇ, ᶜ)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value  = "o-Toluidineᵃ"
$ws.Range("A6").Value  = "DNOPᵇ"
$ws.Range("A8").Value  = "MEHPᵇ"
$ws.Range("A10").Value = "Fluorantheneᵇ"
$ws.Range("A11").Value = "MenthoneNAᶜ"
$ws.Range("A13").Value = "Methoxychlorᵃᵇ"
$ws.Range("A14").Value = "TTBNPᵇ"
$ws.Range("A16").Value = "Benz(a)anthraceneᵃᵇ"
$ws.Range("A17").Value = "N-MeFOSAAᵃ"
$ws.Range("A18").Value = "Terbuthylazineᵇ"
$ws.Range("A19").Value = "MDAᵃᵇᶜ"
$ws.Range("A21").Value = "5-NOTᵃ"
$ws.Range("A22").Value = "Dibutyl phthalateᵇᶜ"
$ws.Range("A24").Value = "o-aminoazotolueneᵃᵇᶜ"
$ws.Range("A25").Value = "OD-PABAᵇ"
$ws.Range("A26").Value = "Prosulfuronᵇ"
$ws.Range("A29").Value = "Metalaxylᵇ"
$ws.Range("A31").Value = "Anthraceneᵃ"

$ws.Range("A33").Value = "ᵃ Possible, likely, or known carcinogen`nᵇ Potential endocrine disrupting chemical`nᶜ Indicates level 2 identification`nAbbreviations: 5-NOT = 5-Nitro-o-toluidine; DEET = N,N-Diethyl-meta-toluamide; DNOP = Di-n-octyl phthalate; MDA = 4,4'-Diaminodiphenylmethane; MEHP = Mono-2-ethylhexyl phthalate; N-MeFOSAA = N-Methylperfluoro-1-octanesulfonamidoacetic acid (linear); OD-PABA = Octyl-dimethyl-p-aminobenzoic acid; PAH = polycyclic aromatic hydrocarbon; TEEP = Tetraethyl ethylenediphosphonate; TTBNP = Tris(tribromoneopentyl); UV = ultraviolet"
